# Add two new columns ("Diploma Institue" / "Diploma Board") right after the
# existing "Diploma %" column (currently AF), pushing every later column two
# places to the right (old AG..AX -> AI..AZ).
#
# The header row (row 1) gets the two new labels; the data rows (2-21) get
# the same "N/A" placeholder text already used for every other column in
# that stretch (AF:AX), keeping the pattern consistent across the newly
# widened AF:AZ block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank columns at AG:AH, shifting everything from the old AG
# column onward to the right by two columns.
$ws.Range("AG1:AH1").EntireColumn.Insert(-4161) | Out-Null  # xlShiftToRight

# New header labels for the inserted columns.
$ws.Range("AG1").Value2 = "Diploma Institue"
$ws.Range("AH1").Value2 = "Diploma Board"

# Give the two new columns a sensible width (matching the sibling columns'
# auto-fit style sizing).
$ws.Columns.Item(33).ColumnWidth = 13.59
$ws.Columns.Item(34).ColumnWidth = 12.25

# Fill the new columns for every data row with the same "N/A" placeholder
# used throughout the rest of that row's AF:AZ span.
for ($r = 2; $r -le 21; $r++) {
  $ws.Cells.Item($r, 33).Value2 = "N/A"
  $ws.Cells.Item($r, 34).Value2 = "N/A"
}

# Restore the last-used selection cell recorded in the saved workbook.
$ws.Range("AV23").Select() | Out-Null
